# edit.ps1 - applies the WoT "Compute Utility" slide edits described in the
# commit diff:
#   1. Slide 1 (title): "Edge Compute Utilities" -> "Compute Utilities"
#      (split into two runs: "Compute " + "Utilities")
#   2. Slide 13 (Summary), bullet: "...performance, etc." ->
#      "...performance, connectivity, etc."
#   3. Slide 13 (Summary), bullet: "Has standardized network interface" ->
#      "...network interface (described by WoT TD, for example)"
#   4. Slide 13 (Summary), bullet: "Has standardized workload packaging" ->
#      "...packaging (using scripts and including WebGPU and WoT Scripting
#      API, for example)" with "WebGPU" split out as its own run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Title slide: "Edge Compute Utilities" -> "Compute Utilities"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# Remove the leading "Edge " (5 characters: E-d-g-e-space)
$edgePrefix = $titleRange.Characters(1, 5)
$edgePrefix.Text = ""

# Re-select "Utilities" (the tail of "Compute Utilities") and re-assign its
# own text so the run is split into "Compute " + "Utilities".
$utilitiesRun = $titleRange.Characters(9, 9)
$utilitiesRun.Text = "Utilities"

# ---------------------------------------------------------------------
# Slide 13: "Summary" content placeholder bullet edits
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$content = $s13.Shapes.Item(3)
$contentRange = $content.TextFrame.TextRange

# 2. "Decision requires metrics on performance, etc." ->
#    "Decision requires metrics on performance, connectivity, etc."
$para4 = $contentRange.Paragraphs(4, 1)
$para4run = $contentRange.Characters($para4.Start, $para4.Length)
$para4run.Text = "Decision requires metrics on performance, connectivity, etc."

# 3. "Has standardized network interface" ->
#    "Has standardized network interface (described by WoT TD, for example)"
$para7 = $contentRange.Paragraphs(7, 1)
$para7run = $contentRange.Characters($para7.Start, $para7.Length)
$para7run.Text = "Has standardized network interface (described by WoT TD, for example)"

# 4. "Has standardized workload packaging" ->
#    "Has standardized workload packaging (using scripts and including
#    WebGPU and WoT Scripting API, for example)" with "WebGPU" broken out
#    into its own run (as PowerPoint's spell-checker would flag it).
$para8 = $contentRange.Paragraphs(8, 1)
$para8Start = $para8.Start
$para8run = $contentRange.Characters($para8Start, $para8.Length)
$newPara8Text = "Has standardized workload packaging (using scripts and including WebGPU and WoT Scripting API, for example)"
$para8run.Text = $newPara8Text

$webGpuOffset = $newPara8Text.IndexOf("WebGPU")
$webGpuStart = $para8Start + $webGpuOffset
$webGpuRun = $contentRange.Characters($webGpuStart, 6)
$webGpuRun.Text = "WebGPU"
